$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 140.14285
$ws.Range("I6").Value = 140.14285
$ws.Range("K6").Value = 420.42855
$ws.Range("M6").Value = -308.42855

$ws.Range("H17").Value = 681.8570999999999
$ws.Range("J17").Value = 681.8570999999999
$ws.Range("L17").Value = 2045.5713
$ws.Range("N17").Value = -2381.5713

$ws.Range("H19").Value = 1953.3846
$ws.Range("I19").Value = 319.8
$ws.Range("J19").Value = 2974.375
$ws.Range("K19").Value = 319.8
$ws.Range("L19").Value = 2974.375
$ws.Range("M19").Value = -144.8
$ws.Range("N19").Value = -3324.375

$ws.Range("H62").Value = 931.3333
$ws.Range("I62").Value = 931.3333
$ws.Range("K62").Value = 931.3333
$ws.Range("M62").Value = -307.3333

$ws.Range("H65").Value = 931.3333
$ws.Range("I65").Value = 931.3333
$ws.Range("K65").Value = 4656.6665
$ws.Range("M65").Value = -1536.6665

$ws.Range("H121").Value = 5646.3335
$ws.Range("J121").Value = 5646.3335
$ws.Range("L121").Value = 16939.0005
$ws.Range("N121").Value = -20433.0005

$ws.Range("H132").Value = 2624.5186
$ws.Range("I132").Value = 1119.75
$ws.Range("K132").Value = 3359.25
$ws.Range("M132").Value = -829.25

$ws.Range("H137").Value = 998.38464
$ws.Range("I137").Value = 998.38464
$ws.Range("K137").Value = 2995.15392
$ws.Range("M137").Value = -445.1539199999997

$ws.Range("H138").Value = 2854.0613
$ws.Range("I138").Value = 1816.9412
$ws.Range("K138").Value = 5450.8236
$ws.Range("M138").Value = -310.8235999999997

$ws.Range("H141").Value = 6399.4
$ws.Range("I141").Value = 4249.25
$ws.Range("K141").Value = 12747.75
$ws.Range("M141").Value = -7567.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 23677.5
$ws.Range("J24").Value = 23677.5
$ws.Range("L24").Value = 23677.5
$ws.Range("N24").Value = -24425.5

$ws.Range("H74").Value = 2930.5
$ws.Range("I74").Value = 2397.5
$ws.Range("K74").Value = 2397.5
$ws.Range("M74").Value = -1523.5

$ws.Range("H77").Value = 2930.5
$ws.Range("I77").Value = 2397.5
$ws.Range("K77").Value = 11987.5
$ws.Range("M77").Value = -7619.5

$ws.Range("H100").Value = 23677.5
$ws.Range("J100").Value = 23677.5
$ws.Range("L100").Value = 23677.5
$ws.Range("N100").Value = -25841.5

$ws.Range("H112").Value = 75000
$ws.Range("J112").Value = 75000
$ws.Range("L112").Value = 75000
$ws.Range("N112").Value = -77954

$ws.Range("H122").Value = 11696
$ws.Range("I122").Value = 7950
$ws.Range("K122").Value = 23850
$ws.Range("M122").Value = -21400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2649.4167
$ws.Range("I94").Value = 1542.8572
$ws.Range("K94").Value = 1542.8572
$ws.Range("M94").Value = -1091.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3598.0344
$ws.Range("I31").Value = 2768.5454
$ws.Range("K31").Value = 2768.5454
$ws.Range("M31").Value = -2473.5454

$ws.Range("H34").Value = 3598.0344
$ws.Range("I34").Value = 2768.5454
$ws.Range("K34").Value = 2768.5454
$ws.Range("M34").Value = -2566.5454

$ws.Range("H43").Value = 11366.125
$ws.Range("J43").Value = 11366.125
$ws.Range("L43").Value = 11366.125
$ws.Range("N43").Value = -11734.125

$ws.Range("H58").Value = 2943.8333
$ws.Range("I58").Value = 1722.625
$ws.Range("K58").Value = 1722.625
$ws.Range("M58").Value = -1519.625

$ws.Range("H99").Value = 11742.143
$ws.Range("I99").Value = 8092.615
$ws.Range("K99").Value = 8092.615
$ws.Range("M99").Value = -6594.615

$ws.Range("H101").Value = 11366.125
$ws.Range("J101").Value = 11366.125
$ws.Range("L101").Value = 11366.125
$ws.Range("N101").Value = -17856.125

$ws.Range("H105").Value = 1086
$ws.Range("I105").Value = 1125.9286
$ws.Range("K105").Value = 1125.9286
$ws.Range("M105").Value = 621.0714

$ws.Range("H126").Value = 11742.143
$ws.Range("I126").Value = 8092.615
$ws.Range("K126").Value = 24277.845
$ws.Range("M126").Value = -21807.845

$ws.Range("H132").Value = 2232.0833
$ws.Range("I132").Value = 2078.5
$ws.Range("K132").Value = 6235.5
$ws.Range("M132").Value = -3705.5

$ws.Range("H136").Value = 2943.8333
$ws.Range("I136").Value = 1722.625
$ws.Range("K136").Value = 5167.875
$ws.Range("M136").Value = -2617.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I5").Value = 4117
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4117
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4005
$ws.Range("N5").ClearContents()

$ws.Range("H80").Value = 3200.5715
$ws.Range("I80").Value = 2351.25
$ws.Range("K80").Value = 2351.25
$ws.Range("M80").Value = -1353.25

$ws.Range("H83").Value = 3200.5715
$ws.Range("I83").Value = 2351.25
$ws.Range("K83").Value = 11756.25
$ws.Range("M83").Value = -6764.25

$ws.Range("H92").Value = 9983
$ws.Range("J92").Value = 9983
$ws.Range("L92").Value = 9983
$ws.Range("N92").Value = -13727

$ws.Range("H122").Value = 87609.914
$ws.Range("I122").Value = 3422.75
$ws.Range("J122").Value = 255984.25
$ws.Range("K122").Value = 10268.25
$ws.Range("L122").Value = 767952.75
$ws.Range("M122").Value = -7818.25
$ws.Range("N122").Value = -772852.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4017.0667
$ws.Range("I7").Value = 3232.5715
$ws.Range("J7").Value = 15000
$ws.Range("K7").Value = 3232.5715
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = -3120.5715
$ws.Range("N7").Value = -15224

$ws.Range("H46").Value = 3294.5454
$ws.Range("I46").Value = 2583.3333
$ws.Range("J46").Value = 4148
$ws.Range("K46").Value = 2583.3333
$ws.Range("L46").Value = 4148
$ws.Range("M46").Value = -2395.3333
$ws.Range("N46").Value = -4524

$ws.Range("H74").Value = 30000
$ws.Range("I74").Value = 30000
$ws.Range("K74").Value = 30000
$ws.Range("M74").Value = -29002

$ws.Range("H77").Value = 30000
$ws.Range("I77").Value = 30000
$ws.Range("K77").Value = 90000
$ws.Range("M77").Value = -85008

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H126").Value = 4017.0667
$ws.Range("I126").Value = 3232.5715
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 9697.7145
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = -7227.7145
$ws.Range("N126").Value = -49940

$ws.Range("H132").Value = 252057.5
$ws.Range("J132").Value = 2222
$ws.Range("L132").Value = 6666
$ws.Range("N132").Value = -11726

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 13199.75
$ws.Range("L105").Value = 13199.75
$ws.Range("N105").Value = -20187.75
